$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs")

# Row 4
$ws.Range("B4").Value = "Homologation"
$ws.Range("C4").Value = "4 x WEBSERVERS 32GB RAM 16vCPU (CONECTA WEB)"
$ws.Range("F4").Value = "c6a.4xlarge"

# Row 5
$ws.Range("C5").Value = "4 x WEBSERVERS 32GB RAM 16vCPU (PORTAL MASSAS)"
$ws.Range("H5").Value = 4

# Row 6
$ws.Range("C6").Value = "2 x WEBSERVERS 16GB RAM 8vCPU (PORTAL WEB)"
$ws.Range("F6").Value = "c6i.2xlarge"
$ws.Range("M6").Value = 230

# Row 7
$ws.Range("C7").Value = "2 x WEBSERVERS 16GB RAM 8vCPU (CONECTA WEB)"
$ws.Range("F7").Value = "c6i.2xlarge"

# Row 8
$ws.Range("B8").Value = "Homologation"
$ws.Range("C8").Value = "2 x WEBSERVERS 32GB RAM 16vCPU (WEBSERVER)"
$ws.Range("D8").Value = "sa-east-1"
$ws.Range("E8").Value = "Windows Server"
$ws.Range("F8").Value = "c6a.4xlarge"
$ws.Range("G8").Value = "Shared Instances"
$ws.Range("H8").Value = 2
$ws.Range("I8").Value = 40
$ws.Range("J8").Value = "Hours/Week"
$ws.Range("K8").Value = "On-Demand"
$ws.Range("L8").Value = "General Purpose SSD (gp3)"
$ws.Range("M8").Value = 230
$ws.Range("P8").Value = "2x Daily"
$ws.Range("Q8").Value = 10

# Row 9
$ws.Range("B9").Value = "Homologation"
$ws.Range("C9").Value = "2 x WEBSERVERS 32GB RAM 16vCPU (TOWER)"
$ws.Range("D9").Value = "sa-east-1"
$ws.Range("E9").Value = "Windows Server"
$ws.Range("F9").Value = "c6a.4xlarge"
$ws.Range("G9").Value = "Shared Instances"
$ws.Range("H9").Value = 2
$ws.Range("I9").Value = 40
$ws.Range("J9").Value = "Hours/Week"
$ws.Range("K9").Value = "On-Demand"
$ws.Range("L9").Value = "General Purpose SSD (gp3)"
$ws.Range("M9").Value = 230
$ws.Range("P9").Value = "2x Daily"
$ws.Range("Q9").Value = 10

# Row 10
$ws.Range("B10").Value = "Homologation"
$ws.Range("C10").Value = "2 x WEBSERVERS 32GB RAM 16vCPU (WEB THINKERS)"
$ws.Range("D10").Value = "sa-east-1"
$ws.Range("E10").Value = "Windows Server"
$ws.Range("F10").Value = "c6a.4xlarge"
$ws.Range("G10").Value = "Shared Instances"
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = 40
$ws.Range("J10").Value = "Hours/Week"
$ws.Range("K10").Value = "On-Demand"
$ws.Range("L10").Value = "General Purpose SSD (gp3)"
$ws.Range("M10").Value = 230
$ws.Range("P10").Value = "2x Daily"
$ws.Range("Q10").Value = 10

# Row 11
$ws.Range("B11").Value = "Homologation"
$ws.Range("C11").Value = "2 x WEBSERVERS 8GB RAM 4vCPU (WEB THINKERS)"
$ws.Range("D11").Value = "sa-east-1"
$ws.Range("E11").Value = "Windows Server"
$ws.Range("F11").Value = "c6in.xlarge"
$ws.Range("G11").Value = "Shared Instances"
$ws.Range("H11").Value = 2
$ws.Range("I11").Value = 40
$ws.Range("J11").Value = "Hours/Week"
$ws.Range("K11").Value = "On-Demand"
$ws.Range("L11").Value = "General Purpose SSD (gp3)"
$ws.Range("M11").Value = 230
$ws.Range("P11").Value = "2x Daily"
$ws.Range("Q11").Value = 10

# Row 12
$ws.Range("B12").Value = "Homologation"
$ws.Range("C12").Value = "2 x WEBSERVERS 32GB RAM 16vCPU (PORTAL GERACAO NEGOCIO)"
$ws.Range("D12").Value = "sa-east-1"
$ws.Range("E12").Value = "Windows Server"
$ws.Range("F12").Value = "c6a.4xlarge"
$ws.Range("G12").Value = "Shared Instances"
$ws.Range("H12").Value = 2
$ws.Range("I12").Value = 40
$ws.Range("J12").Value = "Hours/Week"
$ws.Range("K12").Value = "On-Demand"
$ws.Range("L12").Value = "General Purpose SSD (gp3)"
$ws.Range("M12").Value = 230
$ws.Range("P12").Value = "2x Daily"
$ws.Range("Q12").Value = 10

# Row 13
$ws.Range("B13").Value = "Homologation"
$ws.Range("C13").Value = "2 x BACKUP 32GB RAM 16vCPU (BACKUP)"
$ws.Range("D13").Value = "sa-east-1"
$ws.Range("E13").Value = "Windows Server"
$ws.Range("F13").Value = "c6a.4xlarge"
$ws.Range("G13").Value = "Shared Instances"
$ws.Range("H13").Value = 2
$ws.Range("I13").Value = 40
$ws.Range("J13").Value = "Hours/Week"
$ws.Range("K13").Value = "On-Demand"
$ws.Range("L13").Value = "General Purpose SSD (gp3)"
$ws.Range("M13").Value = 230
$ws.Range("P13").Value = "2x Daily"
$ws.Range("Q13").Value = 10

# Row 14
$ws.Range("B14").Value = "Homologation"
$ws.Range("C14").Value = "2 x WEBSERVERS 16GB RAM 8vCPU (PORTAL GERACAO NEGOCIO)"
$ws.Range("D14").Value = "sa-east-1"
$ws.Range("E14").Value = "Windows Server"
$ws.Range("F14").Value = "c6i.2xlarge"
$ws.Range("G14").Value = "Shared Instances"
$ws.Range("H14").Value = 2
$ws.Range("I14").Value = 40
$ws.Range("J14").Value = "Hours/Week"
$ws.Range("K14").Value = "On-Demand"
$ws.Range("L14").Value = "General Purpose SSD (gp3)"
$ws.Range("M14").Value = 230
$ws.Range("P14").Value = "2x Daily"
$ws.Range("Q14").Value = 10

# Row 15
$ws.Range("B15").Value = "Homologation"
$ws.Range("C15").Value = "2 x WEBSERVERS 16GB RAM 8vCPU (PORTAL MASSAS)"
$ws.Range("D15").Value = "sa-east-1"
$ws.Range("E15").Value = "Windows Server"
$ws.Range("F15").Value = "c6i.2xlarge"
$ws.Range("G15").Value = "Shared Instances"
$ws.Range("H15").Value = 2
$ws.Range("I15").Value = 40
$ws.Range("J15").Value = "Hours/Week"
$ws.Range("K15").Value = "On-Demand"
$ws.Range("L15").Value = "General Purpose SSD (gp3)"
$ws.Range("M15").Value = 170
$ws.Range("P15").Value = "2x Daily"
$ws.Range("Q15").Value = 10

# Row 16
$ws.Range("B16").Value = "Homologation"
$ws.Range("C16").Value = "2 x WEBSERVERS 8GB RAM 4vCPU (PORTAL MASSAS)"
$ws.Range("D16").Value = "sa-east-1"
$ws.Range("E16").Value = "Windows Server"
$ws.Range("F16").Value = "c6in.xlarge"
$ws.Range("G16").Value = "Shared Instances"
$ws.Range("H16").Value = 2
$ws.Range("I16").Value = 40
$ws.Range("J16").Value = "Hours/Week"
$ws.Range("K16").Value = "On-Demand"
$ws.Range("L16").Value = "General Purpose SSD (gp3)"
$ws.Range("M16").Value = 170
$ws.Range("P16").Value = "2x Daily"
$ws.Range("Q16").Value = 10

# Row 17
$ws.Range("B17").Value = "Homologation"
$ws.Range("C17").Value = "2 x WEBSERVERS 16GB RAM 4vCPU (BACKUP)"
$ws.Range("D17").Value = "sa-east-1"
$ws.Range("E17").Value = "Windows Server"
$ws.Range("F17").Value = "m6id.xlarge"
$ws.Range("G17").Value = "Shared Instances"
$ws.Range("H17").Value = 2
$ws.Range("I17").Value = 40
$ws.Range("J17").Value = "Hours/Week"
$ws.Range("K17").Value = "On-Demand"
$ws.Range("L17").Value = "General Purpose SSD (gp3)"
$ws.Range("M17").Value = 230
$ws.Range("P17").Value = "2x Daily"
$ws.Range("Q17").Value = 10

# Row 18
$ws.Range("B18").Value = "Development"
$ws.Range("C18").Value = "2 x WEBSERVERS 32GB RAM 8vCPU (PORTAL QI)"
$ws.Range("D18").Value = "sa-east-1"
$ws.Range("E18").Value = "Windows Server"
$ws.Range("F18").Value = "m6id.2xlarge"
$ws.Range("G18").Value = "Shared Instances"
$ws.Range("H18").Value = 2
$ws.Range("I18").Value = 40
$ws.Range("J18").Value = "Hours/Week"
$ws.Range("K18").Value = "On-Demand"
$ws.Range("L18").Value = "General Purpose SSD (gp3)"
$ws.Range("M18").Value = 170
$ws.Range("P18").Value = "2x Daily"
$ws.Range("Q18").Value = 10

# Row 19
$ws.Range("B19").Value = "Homologation"
$ws.Range("C19").Value = "2 x WEBSERVERS 8GB RAM 4vCPU (PORTAL WEB)"
$ws.Range("D19").Value = "sa-east-1"
$ws.Range("E19").Value = "Windows Server"
$ws.Range("F19").Value = "c6in.xlarge"
$ws.Range("G19").Value = "Shared Instances"
$ws.Range("H19").Value = 2
$ws.Range("I19").Value = 40
$ws.Range("J19").Value = "Hours/Week"
$ws.Range("K19").Value = "On-Demand"
$ws.Range("L19").Value = "General Purpose SSD (gp3)"
$ws.Range("M19").Value = 230
$ws.Range("P19").Value = "2x Daily"
$ws.Range("Q19").Value = 10

# Row 20
$ws.Range("B20").Value = "Development"
$ws.Range("C20").Value = "2 x BACKUP 16GB RAM 8vCPU (PORTAL QI)"
$ws.Range("D20").Value = "sa-east-1"
$ws.Range("E20").Value = "Windows Server"
$ws.Range("F20").Value = "c6i.2xlarge"
$ws.Range("G20").Value = "Shared Instances"
$ws.Range("H20").Value = 2
$ws.Range("I20").Value = 40
$ws.Range("J20").Value = "Hours/Week"
$ws.Range("K20").Value = "On-Demand"
$ws.Range("L20").Value = "General Purpose SSD (gp3)"
$ws.Range("M20").Value = 230
$ws.Range("P20").Value = "2x Daily"
$ws.Range("Q20").Value = 10

